# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (F) and "最低票价" (G) figures for a few events
# across the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 98
$ws1.Range("G2").Value = 55

$ws1.Range("F3").Value = 4074

$ws1.Range("F4").Value = 2377

$ws1.Range("F11").Value = 89

$ws1.Range("F13").Value = 1527

$ws1.Range("F14").Value = 275

$ws1.Range("F15").Value = 2945

$ws1.Range("F16").Value = 202

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 98
$ws4.Range("G2").Value = 55

$ws4.Range("F3").Value = 4074

$ws4.Range("F4").Value = 2377

$ws4.Range("F13").Value = 90

$ws4.Range("F17").Value = 1527

$ws4.Range("F18").Value = 275

$ws4.Range("F19").Value = 2945

$ws4.Range("F20").Value = 202
